# Refresh the NATMI Tnfsf11 -> Tnfrsf11a LR-pair export (Sheet1) with the new TPM-based run.
# A 4th sending/target cluster, "Inflammatory-Mac", is introduced; a new "MuSCs" sending block
# (rows 17-21) is appended, and every ligand/receptor/edge statistic column (E:T) is refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,"ECs","Tnfsf11","Tnfrsf11a","ECs",[double]"1.0",[double]"0.3333333333333333",[double]"0.020131",[double]"0.060393",[double]"0.0058831740909272",[double]"0.006076836619800507",[double]"3.0",[double]"1.0",[double]"2.283376333333333",[double]"6.850129000000001",[double]"0.03598856860670048",[double]"0.03609671873041641",[double]"0.04596664896633333",[double]"0.4136998406970001",[double]"0.0002117270143964963",[double]"0.0002193538622356334"),
    @(3,"ECs","Tnfsf11","Tnfrsf11a","FAPs",[double]"1.0",[double]"0.3333333333333333",[double]"0.020131",[double]"0.060393",[double]"0.0058831740909272",[double]"0.006076836619800507",[double]"2.0",[double]"0.6666666666666666",[double]"0.05457933333333333",[double]"0.163738",[double]"0.0008602314272510669",[double]"0.0008628165296567295",[double]"0.001098736559333333",[double]"0.009888629034",[double]"5.060891245004804e-06",[double]"5.243195083587204e-06"),
    @(4,"ECs","Tnfsf11","Tnfrsf11a","Inflammatory-Mac",[double]"1.0",[double]"0.3333333333333333",[double]"0.020131",[double]"0.060393",[double]"0.0058831740909272",[double]"0.006076836619800507",[double]"3.0",[double]"1.0",[double]"28.334948",[double]"85.004844",[double]"0.4465905182509515",[double]"0.4479325782902666",[double]"0.570410838188",[double]"5.133697543692",[double]"0.002627369766227749",[double]"0.002722013094955949"),
    @(5,"ECs","Tnfsf11","Tnfrsf11a","MuSCs",[double]"1.0",[double]"0.3333333333333333",[double]"0.020131",[double]"0.060393",[double]"0.0058831740909272",[double]"0.006076836619800507",[double]"2.0",[double]"1.0",[double]"0.570287",[double]"1.140574",[double]"0.008988361894356763",[double]"0.006010248693013806",[double]"0.011480447597",[double]"0.068882685582",[double]"5.288009781675704e-05",[double]"3.652329935181444e-05"),
    @(6,"ECs","Tnfsf11","Tnfrsf11a","Resolving-Mac",[double]"1.0",[double]"0.3333333333333333",[double]"0.020131",[double]"0.060393",[double]"0.0058831740909272",[double]"0.006076836619800507",[double]"3.0",[double]"1.0",[double]"32.20407666666667",[double]"96.61223000000001",[double]"0.5075723198207401",[double]"0.5090976377566464",[double]"0.6483002673766668",[double]"5.834702406390001",[double]"0.002986136321241193",[double]"0.003093703168173522"),
    @(7,"FAPs","Tnfsf11","Tnfrsf11a","ECs",[double]"3.0",[double]"1.0",[double]"3.059758333333333",[double]"9.179275",[double]"0.8941975535822989",[double]"0.9236327796800838",[double]"3.0",[double]"1.0",[double]"2.283376333333333",[double]"6.850129000000001",[double]"0.03598856860670048",[double]"0.03609671873041641",[double]"6.986579764052778",[double]"62.87921787647501",[double]"0.03218089000504029",[double]"0.03334011265830466"),
    @(8,"FAPs","Tnfsf11","Tnfrsf11a","FAPs",[double]"3.0",[double]"1.0",[double]"3.059758333333333",[double]"9.179275",[double]"0.8941975535822989",[double]"0.9236327796800838",[double]"2.0",[double]"0.6666666666666666",[double]"0.05457933333333333",[double]"0.163738",[double]"0.0008602314272510669",[double]"0.0008628165296567295",[double]"0.1669995699944445",[double]"1.50299612995",[double]"0.0007692168377625134",[double]"0.0007969256296407685"),
    @(9,"FAPs","Tnfsf11","Tnfrsf11a","Inflammatory-Mac",[double]"3.0",[double]"1.0",[double]"3.059758333333333",[double]"9.179275",[double]"0.8941975535822989",[double]"0.9236327796800838",[double]"3.0",[double]"1.0",[double]"28.334948",[double]"85.004844",[double]"0.4465905182509515",[double]"0.4479325782902666",[double]"86.69809326756668",[double]"780.2828394081001",[double]"0.3993401488730518",[double]"0.4137252123955057"),
    @(10,"FAPs","Tnfsf11","Tnfrsf11a","MuSCs",[double]"3.0",[double]"1.0",[double]"3.059758333333333",[double]"9.179275",[double]"0.8941975535822989",[double]"0.9236327796800838",[double]"2.0",[double]"1.0",[double]"0.570287",[double]"1.140574",[double]"0.008988361894356763",[double]"0.006010248693013806",[double]"1.744940400641667",[double]"10.46964240385",[double]"0.008037371216646175",[double]"0.005551262706896933"),
    @(11,"FAPs","Tnfsf11","Tnfrsf11a","Resolving-Mac",[double]"3.0",[double]"1.0",[double]"3.059758333333333",[double]"9.179275",[double]"0.8941975535822989",[double]"0.9236327796800838",[double]"3.0",[double]"1.0",[double]"32.20407666666667",[double]"96.61223000000001",[double]"0.5075723198207401",[double]"0.5090976377566464",[double]"98.53669194813891",[double]"886.8302275332502",[double]"0.453869926649798",[double]"0.4702192662897357"),
    @(12,"Inflammatory-Mac","Tnfsf11","Tnfrsf11a","ECs",[double]"1.0",[double]"0.3333333333333333",[double]"0.014756",[double]"0.044268",[double]"0.004312359887026068",[double]"0.004454314299427563",[double]"3.0",[double]"1.0",[double]"2.283376333333333",[double]"6.850129000000001",[double]"0.03598856860670048",[double]"0.03609671873041641",[double]"0.03369350117466667",[double]"0.3032415105720001",[double]"0.0001551956596510208",[double]"0.0001607861304033086"),
    @(13,"Inflammatory-Mac","Tnfsf11","Tnfrsf11a","FAPs",[double]"1.0",[double]"0.3333333333333333",[double]"0.014756",[double]"0.044268",[double]"0.004312359887026068",[double]"0.004454314299427563",[double]"2.0",[double]"0.6666666666666666",[double]"0.05457933333333333",[double]"0.163738",[double]"0.0008602314272510669",[double]"0.0008628165296567295",[double]"0.0008053726426666667",[double]"0.007248353784",[double]"3.709627500436684e-06",[double]"3.843256005832437e-06"),
    @(14,"Inflammatory-Mac","Tnfsf11","Tnfrsf11a","Inflammatory-Mac",[double]"1.0",[double]"0.3333333333333333",[double]"0.014756",[double]"0.044268",[double]"0.004312359887026068",[double]"0.004454314299427563",[double]"3.0",[double]"1.0",[double]"28.334948",[double]"85.004844",[double]"0.4465905182509515",[double]"0.4479325782902666",[double]"0.418110492688",[double]"3.762994434192001",[double]"0.001925859036831586",[double]"0.001995232488657791"),
    @(15,"Inflammatory-Mac","Tnfsf11","Tnfrsf11a","MuSCs",[double]"1.0",[double]"0.3333333333333333",[double]"0.014756",[double]"0.044268",[double]"0.004312359887026068",[double]"0.004454314299427563",[double]"2.0",[double]"1.0",[double]"0.570287",[double]"1.140574",[double]"0.008988361894356763",[double]"0.006010248693013806",[double]"0.008415154972",[double]"0.050490929832",[double]"3.876105128329774e-05",[double]"2.677153669640722e-05"),
    @(16,"Inflammatory-Mac","Tnfsf11","Tnfrsf11a","Resolving-Mac",[double]"1.0",[double]"0.3333333333333333",[double]"0.014756",[double]"0.044268",[double]"0.004312359887026068",[double]"0.004454314299427563",[double]"3.0",[double]"1.0",[double]"32.20407666666667",[double]"96.61223000000001",[double]"0.5075723198207401",[double]"0.5090976377566464",[double]"0.4752033552933334",[double]"4.276830197640001",[double]"0.002188834511759726",[double]"0.002267680887664224"),
    @(17,"MuSCs","Tnfsf11","Tnfrsf11a","ECs",[double]"2.0",[double]"1.0",[double]"0.327147",[double]"0.654294",[double]"0.0956069124397477",[double]"0.06583606940068805",[double]"3.0",[double]"1.0",[double]"2.283376333333333",[double]"6.850129000000001",[double]"0.03598856860670048",[double]"0.03609671873041641",[double]"0.7469997173210001",[double]"4.481998303926001",[double]"0.003440755927612666",[double]"0.002376466079472811"),
    @(18,"MuSCs","Tnfsf11","Tnfrsf11a","FAPs",[double]"2.0",[double]"1.0",[double]"0.327147",[double]"0.654294",[double]"0.0956069124397477",[double]"0.06583606940068805",[double]"2.0",[double]"0.6666666666666666",[double]"0.05457933333333333",[double]"0.163738",[double]"0.0008602314272510669",[double]"0.0008628165296567295",[double]"0.017855465162",[double]"0.107132790972",[double]"8.224407074311195e-05",[double]"5.680444892654126e-05"),
    @(19,"MuSCs","Tnfsf11","Tnfrsf11a","Inflammatory-Mac",[double]"2.0",[double]"1.0",[double]"0.327147",[double]"0.654294",[double]"0.0956069124397477",[double]"0.06583606940068805",[double]"3.0",[double]"1.0",[double]"28.334948",[double]"85.004844",[double]"0.4465905182509515",[double]"0.4479325782902666",[double]"9.269693233356001",[double]"55.61815940013601",[double]"0.04269714057484027",[double]"0.02949012031114712"),
    @(20,"MuSCs","Tnfsf11","Tnfrsf11a","MuSCs",[double]"2.0",[double]"1.0",[double]"0.327147",[double]"0.654294",[double]"0.0956069124397477",[double]"0.06583606940068805",[double]"2.0",[double]"1.0",[double]"0.570287",[double]"1.140574",[double]"0.008988361894356763",[double]"0.006010248693013806",[double]"0.186567681189",[double]"0.746270724756",[double]"0.0008593495286105318",[double]"0.0003956911500686516"),
    @(21,"MuSCs","Tnfsf11","Tnfrsf11a","Resolving-Mac",[double]"2.0",[double]"1.0",[double]"0.327147",[double]"0.654294",[double]"0.0956069124397477",[double]"0.06583606940068805",[double]"3.0",[double]"1.0",[double]"32.20407666666667",[double]"96.61223000000001",[double]"0.5075723198207401",[double]"0.5090976377566464",[double]"10.53546706927",[double]"63.21280241562001",[double]"0.04852742233794112",[double]"0.03351698741107292")
)

foreach ($r in $data) {
    $rowNum = $r[0]
    for ($c = 1; $c -lt $r.Length; $c++) {
        $ws.Cells.Item($rowNum, $c).Value = $r[$c]
    }
}
